# Re-colour the deck's live theme (ppt/theme/theme1.xml, the theme bound to
# the slide master / all slides) from the "Integral" (Red Violet) palette to
# the stock "Office Theme" palette.
#
# PowerPoint's object model doesn't expose an "import this whole theme part"
# verb for a COM/automation caller, so - exactly as Microsoft itself
# recommends - we drive the change through the 12-slot
# ThemeColorScheme.Item(i).RGB surface (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), which is what PowerPoint's own Design > Colors picker
# ultimately writes into the theme's <a:clrScheme>.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (the standard default PowerPoint "Office" set)
# expressed as OLE RGB() integers (R + G*256 + B*65536).
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

foreach ($i in 1..12) {
    $tcs.Item($i).RGB = $officeColors[$i]
}
